# Inserts a new weekly price record as row 175 in the "Zanahoria" sheet,
# pushing the existing rows 175-304 down to 176-305 (dimension grows from
# A1:R304 to A1:R305). The new row carries the same fixed/static metadata
# as the record that used to sit at row 175 (market, region, product,
# quality, unit, origin, etc.) but with its own date and price figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 175:304 down to 176:305, duplicating formatting from row 175.
$ws.Rows.Item(175).Insert()

# Populate the newly inserted row 175 with the new weekly record.
$ws.Range("A175").Value = 7
$ws.Range("B175").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C175").Value = "Ñuble"
$ws.Range("D175").Value = 44762
$ws.Range("E175").Value = 16
$ws.Range("F175").Value = 100114013
$ws.Range("G175").Value = "Zanahoria"
$ws.Range("H175").Value = "Sin especificar"
$ws.Range("I175").Value = "Primera"
$ws.Range("J175").Value = 100
$ws.Range("K175").Value = 8000
$ws.Range("L175").Value = 8500
$ws.Range("M175").Value = 8250
$ws.Range("N175").Value = "`$/saco 20 kilos"
$ws.Range("O175").Value = "Provincia de Diguillín"
$ws.Range("P175").Value = 412
$ws.Range("Q175").Value = 20
$ws.Range("R175").Value = "Hortaliza"
